$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 172; this shifts the existing rows 172-188 down to 173-189,
# preserving all of their values (matches the diff, which shows every old row's
# data reappearing one row lower, down through the new row 189).
$ws.Rows("172:172").Insert()

# Populate the newly inserted row 172 with the new weekly price observation.
$ws.Range("A172").Value = 3
$ws.Range("B172").Value = "Femacal de La Calera"
$ws.Range("C172").Value = "Coquimbo"
$ws.Range("D172").Value = 44449
$ws.Range("E172").Value = 5
$ws.Range("F172").Value = 100112040
$ws.Range("G172").Value = "Cilantro"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 160
$ws.Range("K172").Value = 2500
$ws.Range("L172").Value = 2500
$ws.Range("M172").Value = 2500
$ws.Range("N172").Value = "$/docena de atados (3 kilos)"
$ws.Range("O172").Value = "La Cruz"
$ws.Range("P172").Value = 833
$ws.Range("Q172").Value = 3
$ws.Range("R172").Value = "Hortaliza"
